$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Noche y día"
$ws.Range("B4").Value = "Virginia Woolf"
$ws.Range("C4").Value = "Lumen"
